$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Memory Map")

# Bug fix: the start address of each "wide SRAM" / interrupt-controller row
# used to chain off the previous row's *end* address (=C17+1, =C18+1, ...).
# That broke once the row above could legitimately hold a gap, so rework it
# to step by the fixed register width (4 bytes) off the previous row's
# *start* address instead (=B17+4), and carry that same pattern down through
# the remaining rows of the table (19 through 35) so every start/end/hex
# address below recalculates consistently with the revised hardware
# registers sheet.
$ws.Range("B18").Formula = "=B17+4"
for ($r = 19; $r -le 35; $r++) {
    $prev = $r - 1
    $ws.Range("B$r").Formula = "=B$prev+4"
}

# Leave the sheet scrolled back to the top with the cell below the table
# selected, matching where the author finished editing.
$ws.Activate() | Out-Null
$ws.Range("B37").Select() | Out-Null
